# Apply updated odds values to Sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Al-Akhdoud vs Al-Kholood Club)
$ws.Range("F2").Value = 2.86
$ws.Range("G2").Value = 3.35
$ws.Range("H2").Value = 2.38
$ws.Range("I2").Value = 2.74
$ws.Range("J2").Value = 3.25
$ws.Range("K2").Value = 3.85
$ws.Range("P2").Value = 1.86

# Row 3 (Stuttgart vs Eintracht Frankfurt)
$ws.Range("F3").Value = 1.83
$ws.Range("G3").Value = 1.84
$ws.Range("H3").Value = 4.6
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 1.61
$ws.Range("R3").Value = 1.6
$ws.Range("S3").Value = 2.54
$ws.Range("T3").Value = 1.6
$ws.Range("U3").Value = 2.52
$ws.Range("AK3").Value = 17.5
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 8.6

# Row 4 (Dhamk vs Al-Ittihad)
$ws.Range("Q4").Value = 1.53

# Row 7 (Hamburger SV vs Leverkusen)
$ws.Range("F7").Value = 3.6
$ws.Range("G7").Value = 3.7
$ws.Range("P7").Value = 2.32
$ws.Range("R7").Value = 1.55
$ws.Range("AC7").Value = 9.4
$ws.Range("AE7").Value = 25

# Row 8 (Mainz vs FC Heidenheim)
$ws.Range("F8").Value = 1.72
$ws.Range("G8").Value = 1.73
$ws.Range("I8").Value = 6
$ws.Range("P8").Value = 1.9
$ws.Range("Q8").Value = 2.02
$ws.Range("T8").Value = 1.99
$ws.Range("AC8").Value = 9
$ws.Range("AF8").Value = 9.8
$ws.Range("AH8").Value = 24
$ws.Range("AI8").Value = 95
